# Add a new model column ("test5.rstanarm", no-conifers group) into the
# bb soil-moisture model comparison table by inserting a new column F
# (pushing the existing F:J columns to G:K) and filling in the header /
# model-name cells for the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F; this shifts old F:J -> G:K and carries
# formatting from the column to the left, matching the source workbook.
$ws.Columns("F:F").Insert()

# Header row 1 (group label) and row 2 (model name) for the new column.
$ws.Range("F1").Value = $ws.Range("G1").Value()
$ws.Range("F2").Value = "test5.rstanarm"

# New column J (old column I, previously default width) gets narrowed.
$ws.Columns("J:J").ColumnWidth = 15.15

# Restore the selection to match the author's final cursor position.
$ws.Range("F6").Select()
